$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.9999798272297047
$ws.Range("E2").Value = 0.9999798272297047

$ws.Range("D3").Value = 0.9999999999999998
$ws.Range("E3").Value = 0.9999999999999998

$ws.Range("D4").Value = 0.9984259553661755
$ws.Range("E4").Value = 0.9984259553661755

$ws.Range("D5").Value = 0.00005580738928185382
$ws.Range("E5").Value = 0.00005580738928185382

$ws.Range("D6").Value = 0.00000001642102679696414
$ws.Range("E6").Value = 0.00000001642102679696414

$ws.Range("D7").Value = 0.9999999999288467
$ws.Range("E7").Value = 0.00000000007115330546980658

$ws.Range("D8").Value = 0.000000007318220282822119
$ws.Range("E8").Value = 0.9999999926817797

$ws.Range("D9").Value = 0.8953278219507228
$ws.Range("E9").Value = 0.1046721780492772

$ws.Range("D11").Value = 0.999793436747583
$ws.Range("E11").Value = 0.0002065632524169692
$ws.Range("F11").Value = 7.247763156890869
